$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date formats in column F (Fecha expiracion Licencia).
# Previously the values were stored with the day/month order swapped
# (e.g. "2019-22-01"); now they are corrected to a proper yyyy-mm-dd order
# (e.g. "2019-01-22"). The leading apostrophe keeps the value as text so
# Excel does not reinterpret it as a date serial number.
$ws.Range("F2").Value = "'2019-01-22"
$ws.Range("F3").Value = "'2019-02-22"
$ws.Range("F4").Value = "'2019-03-22"
$ws.Range("F5").Value = "'2019-04-22"
$ws.Range("F6").Value = "'2019-05-22"
$ws.Range("F7").Value = "'2019-06-22"

# Fix the date formats in column G (Fecha de Nacimiento) in the same way.
$ws.Range("G2").Value = "'2001-01-22"
$ws.Range("G3").Value = "'2001-02-22"
$ws.Range("G4").Value = "'2001-02-22"
$ws.Range("G5").Value = "'2001-02-22"
$ws.Range("G6").Value = "'2001-02-22"
$ws.Range("G7").Value = "'2001-02-22"

# Update the selected cell on the sheet.
$ws.Range("F22").Select()

$wb.Save()
